# Generate Report for Archive
#
# The localization status for 70d2e4f3-9044-424f-b315-358cf7ddc94d.md moved
# back from "Ready for handoff" to "In Translation" (e.g. a new handoff
# cycle started for that file), so the generated report needs to reflect
# that on the Overview sheet (both language columns) as well as on each
# per-language status sheet.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E4").Value = "In Translation"
$overview.Range("F4").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C4").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C4").Value = "In Translation"
